# Insert a new data row at row 74 (pushing the existing rows 74-193 down to
# 75-194) and populate it with a new Berenjena / Femacal de La Calera record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 74..193 down by one, creating a blank row 74.
$ws.Rows(74).Insert()

# Fill in the new row 74 with the inserted record's values.
$ws.Range("A74").Value = 3
$ws.Range("B74").Value = "Femacal de La Calera"
$ws.Range("C74").Value = "Coquimbo"
$ws.Range("D74").Value = 44540
$ws.Range("E74").Value = 5
$ws.Range("F74").Value = 100112001
$ws.Range("G74").Value = "Berenjena"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 130
$ws.Range("K74").Value = 7500
$ws.Range("L74").Value = 8000
$ws.Range("M74").Value = 7769
$ws.Range("N74").Value = "`$/caja 60 unidades"
$ws.Range("O74").Value = "Región de Arica y Parinacota"
$ws.Range("P74").Value = 129
$ws.Range("Q74").Value = 60
$ws.Range("R74").Value = "Hortaliza"
